# Auto-generated script to apply odds updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("L3").Value = 1.2
$ws.Range("M3").Value = 4.33
$ws.Range("N3").Value = 1.7
$ws.Range("O3").Value = 2.1

# Row 5
$ws.Range("L5").Value = 1.4
$ws.Range("M5").Value = 2.75

# Row 6
$ws.Range("N6").Value = 1.98
$ws.Range("O6").Value = 1.88

# Row 9
$ws.Range("J9").Value = 1.03
$ws.Range("K9").Value = 17
$ws.Range("N9").Value = 1.53
$ws.Range("O9").Value = 2.4

# Row 13
$ws.Range("T13").Value = 6.6
$ws.Range("U13").Value = 10.75
$ws.Range("X13").Value = 22
$ws.Range("Z13").Value = 7
$ws.Range("AB13").Value = 16.5
$ws.Range("AE13").Value = 7.2

# Row 18
$ws.Range("G18").Value = 1.75
$ws.Range("H18").Value = 3.95
$ws.Range("I18").Value = 3.8
$ws.Range("S18").Value = 2.12
$ws.Range("U18").Value = 9.75
$ws.Range("W18").Value = 14.5
$ws.Range("X18").Value = 12.5
$ws.Range("Z18").Value = 14.5
$ws.Range("AA18").Value = 8
$ws.Range("AB18").Value = 13.5
$ws.Range("AE18").Value = 14.5
$ws.Range("AF18").Value = 23
$ws.Range("AG18").Value = 13
$ws.Range("AH18").Value = 55
$ws.Range("AI18").Value = 30
$ws.Range("AJ18").Value = 32

# Row 21
$ws.Range("N21").Value = 1.95
$ws.Range("O21").Value = 1.9

# Row 22
$ws.Range("N22").Value = 2.25
$ws.Range("O22").Value = 1.62

# Row 25
$ws.Range("L25").Value = 1.29
$ws.Range("M25").Value = 3.5
$ws.Range("N25").Value = 1.95
$ws.Range("O25").Value = 1.9

# Row 30
$ws.Range("G30").Value = 1.72
$ws.Range("H30").Value = 3.3
$ws.Range("I30").Value = 4.4
$ws.Range("P30").Value = 1.4
$ws.Range("Q30").Value = 2.4
$ws.Range("R30").Value = 2.06
$ws.Range("U30").Value = 6.2
$ws.Range("W30").Value = 10.75
$ws.Range("X30").Value = 12.5
$ws.Range("Z30").Value = 7.6
$ws.Range("AA30").Value = 5.7
$ws.Range("AB30").Value = 15
$ws.Range("AC30").Value = 80
$ws.Range("AE30").Value = 8.5
$ws.Range("AF30").Value = 18.5
$ws.Range("AG30").Value = 12.5
$ws.Range("AH30").Value = 55
$ws.Range("AI30").Value = 40

# Row 38
$ws.Range("G38").Value = 1.4
$ws.Range("H38").Value = 4.5
$ws.Range("I38").Value = 6.2
$ws.Range("J38").Value = 1.03
$ws.Range("K38").Value = 9
$ws.Range("L38").Value = 1.19
$ws.Range("M38").Value = 4.15
$ws.Range("N38").Value = 1.57
$ws.Range("O38").Value = 2.25
$ws.Range("P38").Value = 1.31
$ws.Range("Q38").Value = 3.2
$ws.Range("R38").Value = 1.82
$ws.Range("S38").Value = 1.9
$ws.Range("T38").Value = 8
$ws.Range("U38").Value = 7.2
$ws.Range("W38").Value = 9.5
$ws.Range("X38").Value = 11
$ws.Range("Y38").Value = 24
$ws.Range("Z38").Value = 9
$ws.Range("AA38").Value = 9.25
$ws.Range("AB38").Value = 18
$ws.Range("AC38").Value = 75
$ws.Range("AD38").Value = 500
$ws.Range("AE38").Value = 19
$ws.Range("AF38").Value = 40
$ws.Range("AG38").Value = 20
$ws.Range("AH38").Value = 120
$ws.Range("AI38").Value = 65
$ws.Range("AJ38").Value = 55
